# fix meipath and add code mei names for score highlight
#
# Adds the "code" mei names under the mc1/mc2/mc3 groups on row 2 of
# Sheet1 (columns N, S, X), which sit beneath the "mc1:"/"mc2:"/"mc3:"
# headers in row 1 and were previously left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N2").Value = "Ending 1"
$ws.Range("S2").Value = "Ending 2"
$ws.Range("X2").Value = "Ending 3"

# Reflect the author's final selection after making the edit.
[void]$ws.Range("O3").Select()
